$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CountryRow($row, $name, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 09:22"

# --- Polonia (row 31): active/recovered cases updated, totals unchanged ---
$ws.Cells.Item(31, 4).Value = 1133
$ws.Cells.Item(31, 5).Value = 7794

# --- Re-rank the Mexico..Indonesia block (rows 34-41) ---
# Singapur jumps up (new totals push it above Dinamarca), the countries that
# used to sit between Mexico and Singapur shift down one row, and Australia's
# numbers are refreshed too.
Set-CountryRow 35 "Singapur" 8014 1426 768 7235 22 0 11
Set-CountryRow 36 "Dinamarca" 7384 0 4141 2888 84 0 355
Set-CountryRow 37 "Noruega" 7103 25 32 6906 58 0 165
Set-CountryRow 38 "Emiratos Arabes Unidos" 6781 0 1286 5454 1 0 41
Set-CountryRow 39 "Chequia" 6746 0 1298 5262 84 0 186
Set-CountryRow 40 "Australia" 6619 7 4258 2290 49 0 71
# Row 41 (Indonesia) is unchanged.

# --- Re-rank the Azerbaiyan..Bosnia y Herzegovina block (rows 73-77) ---
# Armenia jumps up (new totals push it above Eslovenia), Eslovenia/Lituania
# shift down one row each.
Set-CountryRow 74 "Armenia" 1339 48 580 737 30 2 22
Set-CountryRow 75 "Eslovenia" 1330 0 192 1064 26 0 74
Set-CountryRow 76 "Lituania" 1326 28 242 1047 14 2 37
# Row 77 (Bosnia y Herzegovina) is unchanged.
